{"js": "// Insert a new bold, justified paragraph \"Following are the features\"\n// immediately after the \"Features\" paragraph (and before the trailing\n// empty paragraph), matching the formatting of the surrounding headings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the \"Features\" paragraph (trailing space included in original text).\nlet featuresPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Features\") {\n    featuresPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!featuresPara) {\n  throw new Error('Could not locate the \"Features\" paragraph.');\n}\n\n// Insert the new paragraph right after it.\nconst newPara = featuresPara.insertParagraph(\"Following are the features\", Word.InsertLocation.after);\n\n// Match formatting used throughout the diff: justified, bold, size 12 (24 half-points).\nnewPara.alignment = Word.Alignment.justified;\nnewPara.font.bold = true;\nnewPara.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Insert a new bold, justified paragraph \"Following are the features\"\n# immediately after the \"Features\" paragraph (and before the trailing\n# empty paragraph), matching the formatting of the surrounding headings.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Features\" paragraph via Find.\n$findRange = $d.Content\n$find = $findRange.Find\n$find.ClearFormatting()\n$find.Text = \"Features\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate the 'Features' paragraph.\"\n}\n\n$target = $findRange.Start\n\n# Resolve the paragraph index whose range contains the found text (more\n# reliable than indexing straight off the Find range).\n$count = $d.Paragraphs.Count\n$idx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($target -ge $candidate.Range.Start -and $target -lt $candidate.Range.End) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -eq -1) {\n    throw \"Could not resolve paragraph index for the 'Features' paragraph.\"\n}\n\n$featuresPara = $d.Paragraphs.Item($idx)\n\n# Insert a new empty paragraph right after it, then fill in its text/formatting.\n$featuresPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($idx + 1)\n$newPara.Range.Text = \"Following are the features\"\n$newPara.Range.Font.Bold = $true\n$newPara.Range.Font.Size = 12\n$newPara.Format.Alignment = 3   # wdAlignParagraphJustify -> <w:jc w:val=\"both\"/>\n"}
